# Updated cryptos list on Sun Jul 16 19:16:46 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.366.77"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "1.936.67"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").Value = "'0.7710"
$ws.Range("E5").Value = "  +6.05%  "

$ws.Range("D6").Value = "'245.66"
$ws.Range("E6").Value = "  -2.33%  "

$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'27.87"
$ws.Range("E8").Value = "  -0.45%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.3200"
$ws.Range("E9").Value = "  -3.53%  "

$ws.Range("D10").Value = "'0.07026"
$ws.Range("E10").Value = "  -3.61%  "

$ws.Range("D11").Value = "'0.7817"
$ws.Range("E11").Value = "  -3.58%  "

$ws.Range("D12").Value = "'0.08017"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").Value = "1.934.96"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("D14").Value = "'5.354"
$ws.Range("E14").Value = "  -2.50%  "

$ws.Range("D15").Value = "'94.67"
$ws.Range("E15").Value = "  -0.31%  "

$ws.Range("D16").Value = "'14.46"
$ws.Range("E16").Value = "  -4.71%  "

$ws.Range("D17").Value = "30.354.92"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").Value = "'256.18"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("D19").Value = "'0.000007950"
$ws.Range("E19").Value = "  -4.46%  "

$ws.Range("D20").Value = "'5.771"
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("D21").Value = "2.192.02"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("D24").Value = "'6.718"
$ws.Range("E24").Value = "  -3.82%  "

$ws.Range("D25").Value = "'9.534"
$ws.Range("E25").Value = "  -2.57%  "

$ws.Range("D26").Value = "'164.47"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("E27").Value = "  +2.83%  "

$ws.Range("D28").Value = "'19.08"
$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("D29").Value = "'2.270"
$ws.Range("E29").Value = "  -3.73%  "

$ws.Range("D30").Value = "'1.369"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("D31").Value = "'1.514"
$ws.Range("E31").Value = "  -1.59%  "

$ws.Range("D32").Value = "'4.409"
$ws.Range("E32").Value = "  -0.89%  "

$ws.Range("D33").Value = "'4.123"
$ws.Range("E33").Value = "  -2.21%  "

$ws.Range("D34").Value = "'0.05155"
$ws.Range("E34").Value = "  -2.17%  "

$ws.Range("D35").Value = "'1.278"
$ws.Range("E35").Value = "  +0.37%  "

$ws.Range("D36").Value = "'0.7479"

$ws.Range("D37").Value = "'2.782"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("D38").Value = "'0.01954"
$ws.Range("E38").Value = "  -1.24%  "

$ws.Range("D39").Value = "'2.815"
$ws.Range("E39").Value = "  +0.51%  "

$ws.Range("D40").Value = "'78.56"
$ws.Range("E40").Value = "  -1.16%  "

$ws.Range("D41").Value = "'6.410"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").Value = "'0.4500"
$ws.Range("E42").Value = "  -1.75%  "

$ws.Range("D43").Value = "'1.973"
$ws.Range("E43").Value = "  -3.34%  "

$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").Value = "'0.8344"
$ws.Range("E45").Value = "  -1.40%  "

$ws.Range("D46").Value = "'101.07"
$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("D47").Value = "'9.781"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("D48").Value = "'7.508"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").Value = "'981.66"
$ws.Range("E49").Value = "  +10.40%  "

$ws.Range("D50").Value = "'37.20"
$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("D51").Value = "'0.4153"
$ws.Range("E51").Value = "  -1.66%  "
